$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "authority" column/value added to the user table (H1:H4),
# used to start distinguishing admin users (login/logout, announce).
$ws.Range("H4").Value = "authority"

# Reflect the cursor having moved on to the next row after entry.
$ws.Range("H5").Select()
